$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.991.00'
$ws.Range('E2').Value = '  -0.44%  '

# Row 3
$ws.Range('D3').Value = '1.897.08'
$ws.Range('E3').Value = '  +1.57%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.26%  '

# Row 6
$ws.Range('E6').Value = '  +0.05%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5013'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.23%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3912'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.43%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09311'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.84%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.129'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.50%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.97'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.70%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.351'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.60%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.62%  '

# Row 14
$ws.Range('D14').Value = '1.894.94'
$ws.Range('E14').Value = '  +1.47%  '

# Row 15
$ws.Range('E15').Value = '  +0.12%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.285'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.22%  '

# Row 17
$ws.Range('E17').Value = '  -0.89%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.45'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.60%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06574'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.29%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.77'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.95%  '

# Row 21
$ws.Range('E21').Value = '  +0.03%  '

# Row 22
$ws.Range('E22').Value = '  +1.25%  '

# Row 23
$ws.Range('D23').Value = '28.055.91'
$ws.Range('E23').Value = '  -0.43%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.11%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.319'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.76%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.620'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.46%  '

# Row 27
$ws.Range('D27').Value = '2.113.45'
$ws.Range('E27').Value = '  +1.55%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.08%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '156.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.54%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.51%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.077'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.75%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1062'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.84%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.603'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.36%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.625'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.01%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.565'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.44%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06602'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.07%  '

# Row 37
$ws.Range('E37').Value = '  +1.34%  '

# Row 38
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.290'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.99%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2171'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.31%  '

# Row 40
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.222'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.04%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.988'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.13%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6337'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.90%  '

# Row 43
$ws.Range('E43').Value = '  -0.60%  '

# Row 44
$ws.Range('E44').Value = '  +0.02%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.22'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.62%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5955'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.95%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.710'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.38%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.274'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.46%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.027'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.69%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '123.01'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.82%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.175'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.60%  '
